$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: round Ost/Nord coordinates, drop Starttid/Sluttid
$ws.Range("Q3").Value = 331800
$ws.Range("R3").Value = 6626511
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# Row 4: now holds the record previously in row 6 (Vedtrappmossa / Crossocalyx hellerianus),
# with refreshed Id/coords, no Aktivitet, and dropped Starttid/Sluttid
$ws.Range("A4").Value = 111742278
$ws.Range("B4").Value = 94134
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 53
$ws.Range("F4").Value = "Vedtrappmossa"
$ws.Range("G4").Value = "Crossocalyx hellerianus"
$ws.Range("H4").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 331819
$ws.Range("R4").Value = 6626525
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# Row 5: Id swapped with row 7, coords refreshed, drop Starttid/Sluttid
$ws.Range("A5").Value = 111742299
$ws.Range("Q5").Value = 331808
$ws.Range("R5").Value = 6626504
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

# Row 6: now holds the record previously in row 4 (Thomsons trägnagare / Cacotemnus thomsoni),
# with refreshed Id/coords, Aktivitet added, and dropped Starttid/Sluttid
$ws.Range("A6").Value = 111742281
$ws.Range("B6").Value = 4711
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 100299
$ws.Range("F6").Value = "Thomsons trägnagare"
$ws.Range("G6").Value = "Cacotemnus thomsoni"
$ws.Range("H6").Value = "(Kraatz, 1881)"
$ws.Range("M6").Value = "färska gnagspår"
$ws.Range("Q6").Value = 331822
$ws.Range("R6").Value = 6626518
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()

# Row 7: Id swapped with row 5, coords refreshed, drop Starttid/Sluttid
$ws.Range("A7").Value = 111742269
$ws.Range("Q7").Value = 331780
$ws.Range("R7").Value = 6626525
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
